# Updated the P core dimensions
# The "P" series core rows (96-103) had their Effective magnetic path
# length / cross-sectional area / volume columns collapsed from
# "min/max unit" ranges down to a single representative value with the
# unit re-appended (e.g. "12.5 / 13.4 mm" -> "12.5 mm").
#
# Edits are applied column-by-column (B, then C, then D, then E), each
# top-to-bottom across rows 96-103, which mirrors the order the shared
# strings were (re)created in the original authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: Effective magnetic path length ---
$ws.Range("B96").Value = "12.5 mm"
$ws.Range("B97").Value = "15.9 mm"
$ws.Range("B98").Value = "20 mm"
$ws.Range("B99").Value = "25.9 mm"
$ws.Range("B100").Value = "31.6 mm"
$ws.Range("B101").Value = "37.2 mm"
$ws.Range("B102").Value = "45 mm"
$ws.Range("B103").Value = "52 mm"

# --- Column C: Effective core cross-sectional area ---
$ws.Range("C96").Value = "10 mm^2"
$ws.Range("C97").Value = "15.9 mm^2"
$ws.Range("C98").Value = "25 mm^2"
$ws.Range("C99").Value = "43 mm^2"
$ws.Range("C100").Value = "63 mm^2"
$ws.Range("C101").Value = "93 mm^2"
$ws.Range("C102").Value = "136 mm^2"
$ws.Range("C103").Value = "202 mm^2"

# --- Column D: Minimum core cross-sectional area ---
$ws.Range("D96").Value = "0 mm^2"
$ws.Range("D97").Value = "0 mm^2"
$ws.Range("D98").Value = "20 mm^2"
$ws.Range("D99").Value = "0 mm^2"
$ws.Range("D100").Value = "0 mm^2"
$ws.Range("D101").Value = "76.5 mm^2"
$ws.Range("D102").Value = "0 mm^2"
$ws.Range("D103").Value = "0 mm^2"

# --- Column E: Effective core volume ---
$ws.Range("E96").Value = "125 mm^3"
$ws.Range("E97").Value = "253 mm^3"
$ws.Range("E98").Value = "500 mm^3"
$ws.Range("E99").Value = "1114 mm^3"
$ws.Range("E100").Value = "1990 mm^3"
$ws.Range("E101").Value = "3460 mm^3"
$ws.Range("E102").Value = "6120 mm^3"
$ws.Range("E103").Value = "10500 mm^3"

# --- Row heights: the whole body of the table settles back to the
# single-line height (15) now that none of the rows need to wrap. ---
$ws.Range("A24:A79").EntireRow.RowHeight = 15
$ws.Range("A81:A84").EntireRow.RowHeight = 15
$ws.Range("A86:A94").EntireRow.RowHeight = 15
$ws.Range("A96:A115").EntireRow.RowHeight = 15

# --- View state: selection / scroll position left where the editor's
# cursor ended up after making the change. ---
$ws.Range("J98").Select()
$excel.ActiveWindow.ScrollRow = 83
$excel.ActiveWindow.ScrollColumn = 1
